# Adds 5 new device rows (Finger Print Scanner 32, IRIS Scanner 32, Web Camera 32,
# Document Scanner 32, Printer 32) with their mac_address / serial_num / dspec_id
# values to the bottom of the device_master test-data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# lang_code / cr_by / cr_dtimes are constant for every data row already present
# in the sheet - reuse the same values that row 156 (the last existing row) has.
$langCode  = $ws.Range("G156").Value2
$crBy      = $ws.Range("I156").Value2
$crDtimes  = $ws.Range("J156").Value2

$newRows = @(
    @{ Row = 157; Id = 3000176; Name = "Finger Print Scanner 32";  Mac = "80-75-40-E8-CA-24"; Serial = "BS563Q2230824"; DspecId = 165 },
    @{ Row = 158; Id = 3000177; Name = "IRIS Scanner 32";          Mac = "0E-1A-14-4A-6D-3A"; Serial = "BS563Q2230825"; DspecId = 327 },
    @{ Row = 159; Id = 3000178; Name = "Web Camera 32";            Mac = "65-13-7F-0F-F7-53"; Serial = "BS563Q2230826"; DspecId = 736 },
    @{ Row = 160; Id = 3000179; Name = "Document Scanner 32";      Mac = "73-C4-DE-8E-C9-8D"; Serial = "BS563Q2230827"; DspecId = 801 },
    @{ Row = 161; Id = 3000180; Name = "Printer 32";                Mac = "EC-74-AB-E0-0F-38"; Serial = "BS563Q2230828"; DspecId = 920 }
)

# New shared-string entries must be interned in the same order the workbook
# author typed them in: all five device names first, then all five MAC
# addresses, then all five serial numbers.
foreach ($r in $newRows) {
    $ws.Range("B$($r.Row)").Value = $r.Name
}
foreach ($r in $newRows) {
    $ws.Range("C$($r.Row)").Value = $r.Mac
}
foreach ($r in $newRows) {
    $ws.Range("D$($r.Row)").Value = $r.Serial
}

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.Id
    $ws.Range("F$row").Value = $r.DspecId
    $ws.Range("G$row").Value = $langCode

    $ws.Range("H$row").Value = $true
    $ws.Range("H$row").HorizontalAlignment = -4131

    $ws.Range("I$row").Value = $crBy
    $ws.Range("J$row").Value = $crDtimes
}

# Reflect the view state changes that come with this edit: the user scrolled the
# sheet down and selected the columns to the right of the data table.
$ws.Range("K1:XFD1048576").Select()
